$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy($null, $ws1)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "T2A"
Write-Host "done"
